$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the storage description for the larvae PR samples (rows 31-115,
# column G) from the old BIOCODE bag text to the new Molecular bag text,
# wherever it is currently set to the old value.
$oldValue = "BIOCODE -40C mesh bag (2 bags total)"
$newValue = "Molecular -40C mesh bag (2 bags total)"

for ($r = 31; $r -le 115; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}

# Update the view state left behind by the editor: scrolled down so row 82
# is at the top of the frozen pane, with the cursor resting on G118.
$win = $excel.ActiveWindow
$win.ScrollRow = 82
$win.ScrollColumn = 1
$ws.Range("G118").Select()
